# jobactions.xlsx update: add "No" flag to three existing rows and append
# seven new job-application rows (2025-03-15) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the "No" marker (column E) on three existing rows ----------
$ws.Range("E94").Value = "No"
$ws.Range("E98").Value = "No"
$ws.Range("E100").Value = "No"

# --- Append seven new rows (112-118) -------------------------------------
# Column layout: A=Date  B=Entity  C=Jobname  D=Description  F=Web or Email

$newRows = @(
    @{ Row = 112; Entity = "Amazon";  Jobname = "Principal Data Scientist, Forecasting, ASIN Forecasting"; Description = "demand forecasting"; Url = "https://www.linkedin.com/jobs/view/4137885098/?trackingId=e299vDboTfLyewOezp%2F49Q%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D39f6ef31...df0aa05c%29&midToken=AQFpnZsm4rTQjw&midSig=0b_0l8pijg0XI1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m89l2nn0~g1-null-null&eid=1j75g-m89l2nn0-g1&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWU4NDYxZWUyYjA4OWNkZDU0NDliYTQ4NzZmNzdjMTAwNmU0NzViNTY5NThlOWZiN2EzNjhmYWYzOTFhN2JiZmU2OTkzNjc3ZjdlYjYzNWJjY2FhYTY1YTAsMSwx" }
    @{ Row = 113; Entity = "Amazon";  Jobname = "Sr. Data Scientist, Perceptor (Kumo Analytics), AWS Support"; Description = "business stuff?  They mention forecasting"; Url = "https://www.linkedin.com/jobs/view/4184889780/?trackingId=8RyyZELEcPLpz3qE5hd%2BFA%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D39f6ef31...df0aa05c%29&midToken=AQFpnZsm4rTQjw&midSig=0b_0l8pijg0XI1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m89l2nn0~g1-null-null&eid=1j75g-m89l2nn0-g1&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWU4NDYxZWUyYjA4OWNkZDU0NDliYTQ4NzZmNzdjMTAwNmU0NzViNTY5NThlOWZiN2EzNjhmYWYzOTFhN2JiZmU2OTkzNjc3ZjdlYjYzNWJjY2FhYTY1YTAsMSwx" }
    @{ Row = 114; Entity = "Amazon";  Jobname = "Senior Data Scientist, Last Mile Science"; Description = "logistics?  Mention forecasting"; Url = "https://www.linkedin.com/jobs/view/4184242488/?trackingId=EQII1MqLI2PsUdCmISGePw%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D002c1553...c1ec6cde%29&midToken=AQFpnZsm4rTQjw&midSig=0YR_8CJIUw_HE1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m88528ge~mh-null-null&eid=1j75g-m88528ge-mh&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWViNGYxZmVmYjU4ZmNiZDM0OTllYTQ4NzZmNzdjMTAwNmU0NzViNTY4MWFhYWViNDgzMTNiOGIxMzY1NWE1NTY3ZDA1NzdkNmM0NjcyMTI5OWIyYmM3NzcsMSwx" }
    @{ Row = 115; Entity = "Amazon";  Jobname = "Data Scientist, Topline Forecasting"; Description = "Biz forecasting"; Url = "https://www.linkedin.com/jobs/view/4007093976/?trackingId=GXAgJQcAHrymvJVwxCcwfQ%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D96b02801...3f23e2c7%29&midToken=AQFpnZsm4rTQjw&midSig=3GhlV4E8MsYXE1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m85cbosa~6t-null-null&eid=1j75g-m85cbosa-6t&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWViNDExYWVmYmQ4ZWNkZDY0NzlmYTQ4NzZmNzdjMTAwNmU0NzViNTY5OTgyZDRhMTk5NTRmMWZmZTFlMjE1NzI2M2E1MWU0ODBmYjQ3OTI3ZTg2ODA5NjMsMSwx" }
    @{ Row = 116; Entity = "Amazon";  Jobname = "Sr. Data Scientist, Devices Decision Scienc"; Description = "biz stuff, mention forecasting"; Url = "https://www.linkedin.com/jobs/view/4148960536/?trackingId=zHv7uEAnsa0%2FK6LDHqPWow%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D4987540d...389951e1%29&midToken=AQFpnZsm4rTQjw&midSig=2wsph5SHHgRHE1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m7y5jufr~ei-null-null&eid=1j75g-m7y5jufr-ei&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWViNDUxZmUzYjY4OWNhZDU0NjkwYTQ4NzZmNzdjMTAwNmU0NzViNTZiOWY0OTNhMTkwNTRlNWZiNDVkOTEyYTY5ZGYwNWJkOWEwNzRiYjg5MGQ2OTAxNTIsMSwx" }
    @{ Row = 117; Entity = "Amazon";  Jobname = "Sr. Applied Scientist, Renewable Energy Optimization"; Description = "RES opt but want forecasting experience"; Url = "https://www.amazon.jobs/en/jobs/2913322/sr-applied-scientist-renewable-energy-optimization?cmpid=DA_INAD200785B" }
    @{ Row = 118; Entity = "Strella"; Jobname = "Data Scientist"; Description = "produce transport decisions, I think, Seattle"; Url = "https://www.linkedin.com/jobs/view/4175135972/?refId=ByteString(length%3D16%2Cbytes%3D8a8662a5...2cc15d02)&trackingId=Ib7%2BrwAH6ArYZYJILpsqvg%3D%3D" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the date value+format (s="3", numFmtId 14) from the row above so
    # the new date cell reuses the existing style instead of minting a new one.
    $ws.Range("A" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null
    $ws.Range("A" + $row).Value = 45731

    $ws.Range("B" + $row).Value = $r.Entity
    $ws.Range("C" + $row).Value = $r.Jobname
    $ws.Range("D" + $row).Value = $r.Description
    $ws.Range("F" + $row).Value = $r.Url
}

$excel.CutCopyMode = $false

# --- Restore the view: scroll/selection state as of the edit ------------
$ws.Activate()
$ws.Range("E100").Select()
